$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (Resolving-Mac), which also drops the now-unused
# "Resolving-Mac" shared string and re-indexes subsequent rows.
$ws.Rows.Item(4).Delete()

# Update recomputed TPM-derived values on row 2 (Inflammatory-Mac -> Ccl22/Ackr2/FAPs)
$ws.Range("G2").Value = 0.396835
$ws.Range("H2").Value = 1.190505
$ws.Range("I2").Value = 0.4359077351420346
$ws.Range("J2").Value = 0.4359077351420347
$ws.Range("Q2").Value = 0.1039787067
$ws.Range("R2").Value = 0.9358083602999999
$ws.Range("S2").Value = 0.4359077351420346
$ws.Range("T2").Value = 0.4359077351420347

# Update recomputed TPM-derived values on row 3 (Neutrophils -> Ccl22/Ackr2/FAPs)
$ws.Range("G3").Value = 0.5135296666666667
$ws.Range("H3").Value = 1.540589
$ws.Range("I3").Value = 0.5640922648579653
$ws.Range("J3").Value = 0.5640922648579654
$ws.Range("Q3").Value = 0.13455504326
$ws.Range("R3").Value = 1.21099538934
$ws.Range("S3").Value = 0.5640922648579653
$ws.Range("T3").Value = 0.5640922648579654
